$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "合富中国"
$ws.Cells.Item(2, 2).Value = "华夏幸福"
$ws.Cells.Item(2, 3).Value = "平潭发展"
$ws.Cells.Item(3, 1).Value = "平潭发展"
$ws.Cells.Item(3, 2).Value = "平潭发展"
$ws.Cells.Item(3, 3).Value = "孚日股份"
$ws.Cells.Item(4, 1).Value = "华夏幸福"
$ws.Cells.Item(4, 2).Value = "航天发展"
$ws.Cells.Item(4, 3).Value = "合富中国"
$ws.Cells.Item(5, 1).Value = "孚日股份"
$ws.Cells.Item(5, 2).Value = "永泰能源"
$ws.Cells.Item(5, 3).Value = "华夏幸福"
$ws.Cells.Item(6, 1).Value = "众生药业"
$ws.Cells.Item(6, 2).Value = "孚日股份"
$ws.Cells.Item(6, 3).Value = "安泰集团"
$ws.Cells.Item(7, 1).Value = "安泰集团"
$ws.Cells.Item(7, 2).Value = "合富中国"
$ws.Cells.Item(7, 3).Value = "永泰能源"
$ws.Cells.Item(8, 1).Value = "日出东方"
$ws.Cells.Item(8, 2).Value = "日出东方"
$ws.Cells.Item(8, 3).Value = "海马汽车"
$ws.Cells.Item(9, 1).Value = "人民同泰"
$ws.Cells.Item(9, 2).Value = "安泰集团"
$ws.Cells.Item(9, 3).Value = "航天发展"
$ws.Cells.Item(10, 1).Value = "航天发展"
$ws.Cells.Item(10, 2).Value = "东百集团"
$ws.Cells.Item(10, 3).Value = "东百集团"
$ws.Cells.Item(11, 1).Value = "永泰能源"
$ws.Cells.Item(11, 2).Value = "海南海药"
$ws.Cells.Item(11, 3).Value = "多氟多"
$ws.Cells.Item(12, 1).Value = "盈新发展"
$ws.Cells.Item(12, 2).Value = "众生药业"
$ws.Cells.Item(12, 3).Value = "国晟科技"
$ws.Cells.Item(13, 1).Value = "东百集团"
$ws.Cells.Item(13, 2).Value = "海马汽车"
$ws.Cells.Item(13, 3).Value = "胜利股份"
$ws.Cells.Item(14, 1).Value = "海马汽车"
$ws.Cells.Item(14, 2).Value = "人民同泰"
$ws.Cells.Item(14, 3).Value = "人民同泰"
$ws.Cells.Item(15, 1).Value = "特一药业"
$ws.Cells.Item(15, 2).Value = "龙洲股份"
$ws.Cells.Item(15, 3).Value = "日出东方"
$ws.Cells.Item(16, 1).Value = "国晟科技"
$ws.Cells.Item(16, 2).Value = "盈新发展"
$ws.Cells.Item(16, 3).Value = "三木集团"
$ws.Cells.Item(17, 1).Value = "三木集团"
$ws.Cells.Item(17, 2).Value = "大东方"
$ws.Cells.Item(17, 3).Value = "摩恩电气"
$ws.Cells.Item(18, 1).Value = "胜利股份"
$ws.Cells.Item(18, 2).Value = "三木集团"
$ws.Cells.Item(18, 3).Value = "长城军工"
$ws.Cells.Item(19, 1).Value = "多氟多"
$ws.Cells.Item(19, 2).Value = "康芝药业"
$ws.Cells.Item(19, 3).Value = "盈新发展"
$ws.Cells.Item(20, 1).Value = "龙洲股份"
$ws.Cells.Item(20, 2).Value = "胜利股份"
$ws.Cells.Item(20, 3).Value = "天际股份"
$ws.Cells.Item(21, 1).Value = "先导智能"
$ws.Cells.Item(21, 2).Value = "荣盛发展"
$ws.Cells.Item(21, 3).Value = "龙洲股份"
